$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws3 = $wb.Worksheets.Item(3)

# Rename "Hoja3" -> "producto a buscar"
$ws3.Name = "producto a buscar"

# Header cell: reuse the workbook's existing bold+bordered header style
# (copy format from url!A1, which already carries that exact style) then set the text.
$ws1.Range("A1").Copy()
$ws3.Range("A1").PasteSpecial(-4122)
$ws3.Range("A1").Value = "Producto"

# Value cell: start from the same bordered format, then drop Bold so it
# collapses onto a plain-font/bordered style, then set the text.
$ws3.Range("A1").Copy()
$ws3.Range("A2").PasteSpecial(-4122)
$ws3.Range("A2").Value = "pantalones"
$ws3.Range("A2").Font.Bold = $false

# Page setup for the new sheet
$ws3.PageSetup.PaperSize = 9
$ws3.PageSetup.Orientation = 1

# Make "producto a buscar" the active/selected sheet & cell
$ws3.Activate() | Out-Null
$ws3.Range("B8").Select() | Out-Null
